$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page 1")

# Update the phone numbers in row 4 with new values
$ws.Range("B4").Value = "448-59-26"
$ws.Range("C4").Value = "775-96-85"
$ws.Range("D4").Value = "816-93-57"
$ws.Range("E4").Value = "402-74-22"
$ws.Range("F4").Value = "976-19-75"
$ws.Range("G4").Value = "528-96-42"

# Remove the now-empty row 5 (it only had blank placeholder cells)
$ws.Rows.Item(5).Delete()
